# Generate Report for Archive
# Update localization status for the two files whose translation moved from
# "Ready for handoff" to "In Translation" (rows for
# 1dbbcbcf-85f4-4818-8198-aace138e8828.md and
# 6d729223-c835-4439-9ba7-83e39f6c319c.md).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: Status columns for zh-cn (B) and de-de (C) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("B4").Value = $newStatus
$wsOverview.Range("C4").Value = $newStatus

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("C4").Value = $newStatus

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("C4").Value = $newStatus
